# Updates market-price-derived profit columns (H:N) across the Leve Profits
# workbook's per-job sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) to reflect
# refreshed Universalis market data pulled by the scheduled runner.
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ

$wb = $excel.ActiveWorkbook

# ALC row 40 - Stuck in the Moment (Horn Glue)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2371.2856
$ws.Range("I40").Value = 3386.75
$ws.Range("J40").Value = 1746.3846
$ws.Range("K40").Value = 3386.75
$ws.Range("L40").Value = 1746.3846
$ws.Range("M40").Value = -3211.75
$ws.Range("N40").Value = -2096.3846

# ALC row 101 - Edge of the Arcane (Cunning Craftsman's Tea)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 877
$ws.Range("I101").Value = 300
$ws.Range("J101").Value = 3185
$ws.Range("K101").Value = 900
$ws.Range("L101").Value = 9555
$ws.Range("M101").Value = 722
$ws.Range("N101").Value = -12799

# ALC row 107 - Another Man's Ink (Enchanted Truegold Ink)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 760.7143
$ws.Range("I107").Value = 742.3077
$ws.Range("K107").Value = 742.3077
$ws.Range("M107").Value = 1177.6923

# ALC row 129 - Practical Command (Commanding Craftsman's Draught)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 950.1525
$ws.Range("I129").Value = 342
$ws.Range("J129").Value = 994.38184
$ws.Range("K129").Value = 1026
$ws.Range("L129").Value = 2983.14552
$ws.Range("M129").Value = 3974
$ws.Range("N129").Value = -12983.14552

# ALC row 137 - Cutting Edge of Culinary Quality (Magnesia Whetstone)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 577307.5600000001
$ws.Range("I137").Value = 3032.7827
$ws.Range("J137").Value = 954688.2
$ws.Range("K137").Value = 9098.348100000001
$ws.Range("L137").Value = 2864064.6
$ws.Range("M137").Value = -6548.348100000001
$ws.Range("N137").Value = -2869164.6

# ARM row 43 - They've Got Legs (Steel Sabatons)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 15338.5
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 15338.5
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 15338.5
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -15964.5

# ARM row 45 - Hollow Hallmarks (Mythril Ingot)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2019.8695
$ws.Range("I45").Value = 2025.35
$ws.Range("K45").Value = 2025.35
$ws.Range("M45").Value = -1648.35

# ARM row 61 - Dealing with the Tough Stuff (Cobalt Ingot)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5752.5884
$ws.Range("I61").Value = 2154.8975
$ws.Range("K61").Value = 2154.8975
$ws.Range("M61").Value = -1942.8975

# ARM row 74 - As the Bolt Flies (Titanium Nugget)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4248.171
$ws.Range("I74").Value = 1621.6786
$ws.Range("J74").Value = 9905.23
$ws.Range("K74").Value = 1621.6786
$ws.Range("L74").Value = 9905.23
$ws.Range("M74").Value = -747.6786
$ws.Range("N74").Value = -11653.23

# ARM row 77 - Heavy Metal Banned (L) (Titanium Nugget)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4248.171
$ws.Range("I77").Value = 1621.6786
$ws.Range("J77").Value = 9905.23
$ws.Range("K77").Value = 8108.393
$ws.Range("L77").Value = 49526.14999999999
$ws.Range("M77").Value = -3740.393
$ws.Range("N77").Value = -58262.14999999999

# ARM row 122 - Haste for High Durium (High Durium Nugget)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1831.6666
$ws.Range("I122").Value = 1789.5652
$ws.Range("K122").Value = 5368.6956
$ws.Range("M122").Value = -2918.6956

# ARM row 132 - Don't Bore Me, Ore Me (Mountain Chromite Ingot)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4910.8
$ws.Range("I132").Value = 5555
$ws.Range("J132").Value = 4749.75
$ws.Range("K132").Value = 16665
$ws.Range("L132").Value = 14249.25
$ws.Range("M132").Value = -14135
$ws.Range("N132").Value = -19309.25

# ARM row 136 - Metal with Mettle (Cobalt Tungsten Ingot)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5752.5884
$ws.Range("I136").Value = 2154.8975
$ws.Range("K136").Value = 6464.6925
$ws.Range("M136").Value = -3914.6925

# BSM row 134 - Ruthenium Supremium (Ruthenium Ingot)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 45611
$ws.Range("I134").Value = 2417.7368
$ws.Range("K134").Value = 7253.2104
$ws.Range("M134").Value = -4718.2104

# CRP row 31 - Wall Not Found (Walnut Lumber)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 474240.66
$ws.Range("I31").Value = 7673.0356
$ws.Range("J31").Value = 746405.1
$ws.Range("K31").Value = 7673.0356
$ws.Range("L31").Value = 746405.1
$ws.Range("M31").Value = -7378.0356
$ws.Range("N31").Value = -746995.1

# CRP row 34 - Armoires of the Rich and Famous (Walnut Lumber)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 474240.66
$ws.Range("I34").Value = 7673.0356
$ws.Range("J34").Value = 746405.1
$ws.Range("K34").Value = 7673.0356
$ws.Range("L34").Value = 746405.1
$ws.Range("M34").Value = -7471.0356
$ws.Range("N34").Value = -746809.1

# CRP row 132 - Hull Lotta Damage (Ginseng Lumber)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4843.4736
$ws.Range("I132").Value = 4917.8335
$ws.Range("K132").Value = 14753.5005
$ws.Range("M132").Value = -12223.5005

# CUL row 2 - Pork Is a Salty Food (Table Salt)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 44.526318
$ws.Range("I2").Value = 10.8
$ws.Range("K2").Value = 64.80000000000001
$ws.Range("M2").Value = 48.19999999999999

# CUL row 38 - Pretty as a Picture (Dark Vinegar)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 67.86957
$ws.Range("I38").Value = 27.785715
$ws.Range("J38").Value = 130.22223
$ws.Range("K38").Value = 83.357145
$ws.Range("L38").Value = 390.66669
$ws.Range("M38").Value = 263.642855
$ws.Range("N38").Value = -1084.66669

# CUL row 68 - Such a Butter Face (Fermented Butter)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 119773.74
$ws.Range("I68").Value = 240539.98
$ws.Range("J68").Value = 3936.7346
$ws.Range("K68").Value = 721619.9400000001
$ws.Range("L68").Value = 11810.2038
$ws.Range("M68").Value = -720808.9400000001
$ws.Range("N68").Value = -13432.2038

# CUL row 71 - No Margarine of Error (L) (Fermented Butter)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 119773.74
$ws.Range("I71").Value = 240539.98
$ws.Range("J71").Value = 3936.7346
$ws.Range("K71").Value = 2164859.82
$ws.Range("L71").Value = 35430.61139999999
$ws.Range("M71").Value = -2160803.82
$ws.Range("N71").Value = -43542.61139999999

# GSM row 35 - Necklet of Champions (Horn Necklace)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 10690
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 10690
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 10690
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -11286

# GSM row 102 - Put the Metal to the Peddle (Durium Ingot)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3918
$ws.Range("I102").Value = 3419.6667
$ws.Range("J102").Value = 4665.5
$ws.Range("K102").Value = 3419.6667
$ws.Range("L102").Value = 4665.5
$ws.Range("M102").Value = -1797.6667
$ws.Range("N102").Value = -7909.5

# GSM row 126 - Gold Rush Order (Phrygian Gold Ingot)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2585.16
$ws.Range("I126").Value = 1708
$ws.Range("J126").Value = 3274.3572
$ws.Range("K126").Value = 5124
$ws.Range("L126").Value = 9823.071599999999
$ws.Range("M126").Value = -2654
$ws.Range("N126").Value = -14763.0716

# GSM row 132 - On Board for Lar (Lar Ingot)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 13278
$ws.Range("I132").Value = 6341.3335
$ws.Range("J132").Value = 15590.223
$ws.Range("K132").Value = 19024.0005
$ws.Range("L132").Value = 46770.669
$ws.Range("M132").Value = -16494.0005
$ws.Range("N132").Value = -51830.669

# LTW row 40 - Best Served Toad (Toad Leather)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3363.48
$ws.Range("I40").Value = 2953.8823
$ws.Range("J40").Value = 4233.875
$ws.Range("K40").Value = 2953.8823
$ws.Range("L40").Value = 4233.875
$ws.Range("M40").Value = -2817.8823
$ws.Range("N40").Value = -4505.875

# LTW row 122 - Hell on Leather (Gaja Leather)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6371.5713
$ws.Range("I122").Value = 6113.3335
$ws.Range("J122").Value = 7921
$ws.Range("K122").Value = 18340.0005
$ws.Range("L122").Value = 23763
$ws.Range("M122").Value = -15890.0005
$ws.Range("N122").Value = -28663

# LTW row 132 - Tenets of Tanning (Silver Lobo Leather)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 54669
$ws.Range("I132").Value = 79004
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 237012
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -234482
$ws.Range("N132").Value = -23057

# LTW row 140 - Worqor Zormor or Bust (Gargantuaskin Shoes of Healing)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 66307.57000000001
$ws.Range("J140").Value = 66307.57000000001
$ws.Range("L140").Value = 66307.57000000001
$ws.Range("N140").Value = -76667.57000000001

# WVR row 132 - Comfy Cabins (Snow Cotton Cloth)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4021.0557
$ws.Range("I132").Value = 4423
$ws.Range("J132").Value = 3217.1667
$ws.Range("K132").Value = 13269
$ws.Range("L132").Value = 9651.500100000001
$ws.Range("M132").Value = -10739
$ws.Range("N132").Value = -14711.5001

# WVR row 136 - Weaving the Envelope (Sarcenet Cloth)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 6134.9697
$ws.Range("I136").Value = 6083.0605
$ws.Range("J136").Value = 6186.879
$ws.Range("K136").Value = 18249.1815
$ws.Range("L136").Value = 18560.637
$ws.Range("M136").Value = -15699.1815
$ws.Range("N136").Value = -23660.637

Write-Output "Updated 190 cells across 31 rows in 8 sheets."